# Regenerate the "K" column (column G) values in the save_data sheet.
# The workbook's strikeout-based column was recalculated ("K" instead of
# "Strike#"), so the raw per-game values in column G change while the
# rest of the row data (TB, PC, dS0, dSF, IP, I0, IF, etc.) stays intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newKValues = @{
    2  = 4
    3  = 2
    4  = 6
    5  = 3
    6  = 8
    7  = 3
    8  = 1
    9  = 6
    10 = 6
    11 = 5
    12 = 5
    13 = 1
    14 = 4
    15 = 3
    16 = 4
    17 = 1
    18 = 2
    19 = 2
    20 = 3
    21 = 2
    22 = 2
    23 = 0
    24 = 1
    25 = 0
    26 = 0
    27 = 1
    28 = 4
    29 = 2
    30 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 2
    36 = 1
    38 = 1
}

foreach ($row in $newKValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newKValues[$row]
}
